{"js": "// Apply the three text corrections to the \"Areas covered\" / skills section\n// of the CV, as described by the diff:\n//  1. Frontend skills list text update\n//  2. Backend skills list text update (append \"...\")\n//  3. Tools line text update\n\nconst replacements = [\n  {\n    oldText:\n      \"JavaScript, CSS, HTML, React, React Router, Redux/Redux toolkit, SASS/SCSS, Bootstrap, Tailwind CSS, DaisyUI\",\n    newText:\n      \"JavaScript, CSS, HTML, React, React Router, Redux/Redux, Toolkit, Sass/SCSS, Bootstrap, Tailwind CSS, DaisyUI, Axios...\"\n  },\n  {\n    oldText:\n      \"Node.js, Express, MySQL, MongoDB, Mongoose, JSON Web Token, Bcrypt\",\n    newText:\n      \"Node.js, Express, MySQL, MongoDB, Mongoose, JSON Web Token, Bcrypt...\"\n  },\n  {\n    oldText:\n      \"(Tools include: GitHub/GitHub desktop, Postman, MySQL Workshop, XAMPP, MongoDB compass)\",\n    newText:\n      \"(Tools include: GitHub/GitHub Desktop, Postman, MySQL Workbench, XAMPP, MongoDB Compass, Firebase/Firestore, PayPal Sandbox, Heroku...)\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the three text corrections to the \"Areas covered\" / skills section\n# of the CV, as described by the diff:\n#  1. Frontend skills list text update\n#  2. Backend skills list text update (append \"...\")\n#  3. Tools line text update\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n\nReplace-Text `\n    \"JavaScript, CSS, HTML, React, React Router, Redux/Redux toolkit, SASS/SCSS, Bootstrap, Tailwind CSS, DaisyUI\" `\n    \"JavaScript, CSS, HTML, React, React Router, Redux/Redux, Toolkit, Sass/SCSS, Bootstrap, Tailwind CSS, DaisyUI, Axios...\"\n\nReplace-Text `\n    \"Node.js, Express, MySQL, MongoDB, Mongoose, JSON Web Token, Bcrypt\" `\n    \"Node.js, Express, MySQL, MongoDB, Mongoose, JSON Web Token, Bcrypt...\"\n\nReplace-Text `\n    \"(Tools include: GitHub/GitHub desktop, Postman, MySQL Workshop, XAMPP, MongoDB compass)\" `\n    \"(Tools include: GitHub/GitHub Desktop, Postman, MySQL Workbench, XAMPP, MongoDB Compass, Firebase/Firestore, PayPal Sandbox, Heroku...)\"\n"}
